$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47 will end up with the same "present but empty" F/G cells that row 45
# currently has (placeholders for a repair time/duration that was never
# filled in). Copy that exact blank pattern down to F47:G47 before we clear
# it out of row 45.
$ws.Range("F45:G45").Copy($ws.Range("F47"))

# --- Row 45: this incident (WC48 P5F / Cámara no detecta busbar) never
# got its "Hora de Reparación" / "Tiempo de Reparación" placeholders filled
# in, and they should no longer be present at all in the saved row. ---
$ws.Range("F45").ClearContents()
$ws.Range("G45").ClearContents()

# --- Row 46: new incident (WV50 FILTER / Traza) ---
$ws.Range("A46").Value = "WV50 FILTER"
$ws.Range("B46").Value = "Traza"
$ws.Range("C46").Value = "'2024-05-30"
$ws.Range("D46").Value = "18:39:20"
$ws.Range("E46").Value = "Noche"
$ws.Range("H46").Value = "N/A"

# --- Row 47: new incident (WC48 P5F / Cámara no detecta Pcb) ---
$ws.Range("A47").Value = "WC48 P5F"
$ws.Range("B47").Value = "Cámara no detecta Pcb"
$ws.Range("C47").Value = "'2024-05-30"
$ws.Range("D47").Value = "18:43:02"
$ws.Range("E47").Value = "Noche"
$ws.Range("H47").Value = "N/A"
